$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry per-row observation data and need to shift down by one
# row (row r's new data = row (r-1)'s old data), for rows 4..25. Row 2 is
# the anchor (unchanged) and row 3 receives a brand-new observation. Row 26
# is a newly appended row that takes the old row 25's data.
$cols = @(4, 10, 11, 12, 13, 15, 16)   # D, J, K, L, M, O, P

# Capture the "before" values for rows 2..25 for each of those columns.
$old = @{}
for ($r = 2; $r -le 25; $r++) {
    $old[$r] = @{}
    foreach ($c in $cols) {
        $old[$r][$c] = $ws.Cells.Item($r, $c).Value2
    }
}

# Row 3 gets brand-new data (not derived from any existing row).
$new3 = @{4 = 44545; 10 = 140; 11 = 14000; 12 = 15000; 13 = 14429; 15 = "Provincia de Chacabuco"; 16 = 577}

# Write row 3's new values.
foreach ($c in $cols) {
    $ws.Cells.Item(3, $c).Value = $new3[$c]
}

# Rows 4..25 take on what used to be in row (r-1).
for ($r = 25; $r -ge 4; $r--) {
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $old[$r - 1][$c]
    }
}

# Row 26 is newly appended, taking the old row 25's data.
foreach ($c in $cols) {
    $ws.Cells.Item(26, $c).Value = $old[25][$c]
}
# Match the date-time number format used by the other rows in column D.
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Fill in the rest of row 26's (constant-across-rows) columns, mirroring
# row 25's structure: A, B, C, E, F, G, H, I, N, Q, R.
$ws.Cells.Item(26, 1).Value = $ws.Cells.Item(25, 1).Value2
$ws.Cells.Item(26, 2).Value = $ws.Cells.Item(25, 2).Value2
$ws.Cells.Item(26, 3).Value = $ws.Cells.Item(25, 3).Value2
$ws.Cells.Item(26, 5).Value = $ws.Cells.Item(25, 5).Value2
$ws.Cells.Item(26, 6).Value = $ws.Cells.Item(25, 6).Value2
$ws.Cells.Item(26, 7).Value = $ws.Cells.Item(25, 7).Value2
$ws.Cells.Item(26, 8).Value = $ws.Cells.Item(25, 8).Value2
$ws.Cells.Item(26, 9).Value = $ws.Cells.Item(25, 9).Value2
$ws.Cells.Item(26, 14).Value = $ws.Cells.Item(25, 14).Value2
$ws.Cells.Item(26, 17).Value = $ws.Cells.Item(25, 17).Value2
$ws.Cells.Item(26, 18).Value = $ws.Cells.Item(25, 18).Value2
